$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("The Order", "Crime", "'7", "'2024")
    ,@("Bird", "Drama", "'7.2", "'2024")
    ,@("The Lord of the Rings: The War of the Rohirrim", "Action", "'6.6", "'2024")
    ,@("Gladiator II", "Action", "'6.8", "'2024")
    ,@("[DA] The Promised Land", "Action", "'7.7", "'2023")
    ,@("The Six Triple Eight", "Drama", "'6.5", "'2024")
    ,@("DragonHeart", "Action", "'6.4", "'1996")
    ,@("Megalopolis", "Drama", "'4.8", "'2024")
    ,@("Hush", "Action", "'6.6", "'2016")
    ,@("Candyman", "Action", "'6.7", "'1992")
    ,@("Teenage Mutant Ninja Turtles: Out of the Shadows", "Action", "'5.9", "'2016")
    ,@("Demolition Man", "Action", "'6.7", "'1993")
    ,@("[ES] The Settlers", "Action", "'6.9", "'2023")
    ,@("Twin Peaks: Fire Walk with Me", "Action", "'7.3", "'1992")
    ,@("The Last House on the Left", "Action", "'6.5", "'2009")
    ,@("Ted", "Action", "'6.9", "'2012")
    ,@("Starve Acre", "Action", "'5.4", "'2023")
    ,@("The Forge", "Action", "'6.4", "'2024")
    ,@("Sideways", "Action", "'7.5", "'2004")
    ,@("[TH] How to Make Millions Before Grandma Dies", "Comedy", "'8", "'2024")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
}
